$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 2 values
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 26.88718205824783
$ws.Range("G2").Value = 20.9173624646147
$ws.Range("H2").Value = 32.34948552140182
$ws.Range("I2").Value = 0.0007645371881185627
$ws.Range("J2").Value = 0.0006267610152500065
$ws.Range("K2").Value = 0.0009715678201341332
$ws.Range("L2").Value = 0.04904985511492113
$ws.Range("M2").Value = 0.04258367460904983
$ws.Range("N2").Value = 0.05428089295621556

# Update Row 3 values
$ws.Range("F3").Value = 0.2086306804287965
$ws.Range("G3").Value = 0.0008641253673257766
$ws.Range("H3").Value = 0.541635818925212
$ws.Range("I3").Value = 0.1934749701365729
$ws.Range("J3").Value = 0.0007903867403991788
$ws.Range("K3").Value = 0.5032470635380228
$ws.Range("L3").Value = 0.2191109099153959
$ws.Range("M3").Value = 0.000923110736303278
$ws.Range("N3").Value = 0.5675610581032168

# Add new Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 27.09581273867663
$ws.Range("G4").Value = 20.91822658998202
$ws.Range("H4").Value = 32.89112134032703
$ws.Range("I4").Value = 0.1942395073246914
$ws.Range("J4").Value = 0.001417147755649185
$ws.Range("K4").Value = 0.504218631358157
$ws.Range("L4").Value = 0.268160765030317
$ws.Range("M4").Value = 0.0435067853453531
$ws.Range("N4").Value = 0.6218419510594324

# Copy style from A3 to A4 (border + bold + alignment)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
